$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data (scraped via a JSON feed) got re-run and appended three
# more rows (duplicates of the existing rows 3, 2 and 4, in that order) to
# the bottom of the sheet, extending the used range from A1:K4 to A1:K7.
# All values in this sheet -- including the numeric-looking ones -- are
# stored as text, so force a text number format on the target range before
# writing the values in order to avoid Excel auto-converting them to
# numbers.

# KL Rahul's name in this sheet contains a non-breaking space (U+00A0)
# between "Rahul" and "(c)", matching the existing rows.
$batsman = "KL Rahul$([char]0x00A0)(c)"

$newRows = @(
    @{ Row = 5; Venue = " Dubai (DSC)";   Date = " October 04 2020"; Result = "Super Kings won by 10 wickets (with 14 balls remaining)"; OwnTeam = "Kings XI Punjab"; OppTeam = "Chennai Super Kings";   Batsman = $batsman; Runs = "63"; Balls = "52"; Fours = "7"; Sixes = "1"; Sr = "121.15" },
    @{ Row = 6; Venue = " Abu Dhabi";     Date = " October 10 2020"; Result = "KKR won by 2 runs";                                        OwnTeam = "Kings XI Punjab"; OppTeam = "Kolkata Knight Riders"; Batsman = $batsman; Runs = "74"; Balls = "58"; Fours = "6"; Sixes = "0"; Sr = "127.58" },
    @{ Row = 7; Venue = " Dubai (DSC)";   Date = " October 08 2020"; Result = "Sunrisers won by 69 runs";                                 OwnTeam = "Kings XI Punjab"; OppTeam = "Sunrisers Hyderabad"; Batsman = $batsman; Runs = "11"; Balls = "16"; Fours = "0"; Sixes = "0"; Sr = "68.75" }
)

# Force every new cell to be treated as plain text (matches the rest of
# the sheet, where even the numeric columns are stored as text strings).
$ws.Range("A5:K7").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.Venue
    $ws.Cells.Item($row, 2).Value  = $r.Date
    $ws.Cells.Item($row, 3).Value  = $r.Result
    $ws.Cells.Item($row, 4).Value  = $r.OwnTeam
    $ws.Cells.Item($row, 5).Value  = $r.OppTeam
    $ws.Cells.Item($row, 6).Value  = $r.Batsman
    $ws.Cells.Item($row, 7).Value  = $r.Runs
    $ws.Cells.Item($row, 8).Value  = $r.Balls
    $ws.Cells.Item($row, 9).Value  = $r.Fours
    $ws.Cells.Item($row, 10).Value = $r.Sixes
    $ws.Cells.Item($row, 11).Value = $r.Sr
}
